$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match results synced from the tracker data source (rows 23-34)
$newRows = @(
    @('23', '14339214', '2025-08-04', 'Clement Chidekh', 'Harold Mayot', 'Gana Clement Chidekh', 2.1),
    @('24', '14339495', '2025-08-04', 'Carlos Taberner', 'Federico Bondioli', 'Gana Federico Bondioli', 5),
    @('25', '14339492', '2025-08-04', 'Gabriele Piraino', 'Jelle Sels', 'Gana Jelle Sels', 2.2),
    @('26', '14339490', '2025-08-04', 'Kimmer Coppejans', 'Tiago Pereira', 'Gana Tiago Pereira', 3.5),
    @('27', '14339485', '2025-08-04', 'Murkel Dellien', 'Dusan Lajovic', 'Gana Murkel Dellien', 3.4),
    @('28', '14339491', '2025-08-04', 'Oleg Prihodko', 'Stefano Travaglia', 'Gana Oleg Prihodko', 2.1),
    @('29', '14339487', '2025-08-04', 'Santiago Rodriguez Taverna', 'Nikolas Sanchez Izquierdo', 'Gana Santiago Rodriguez Taverna', 1.83),
    @('30', '14339504', '2025-08-04', 'Benjamin Hassan', 'Filip Cristian Jianu', 'Gana Filip Cristian Jianu', 2.75),
    @('31', '14339502', '2025-08-04', 'Jan Choinski', 'Geoffrey Blancaneaux', 'Gana Geoffrey Blancaneaux', 2.63),
    @('32', '14339501', '2025-08-04', 'Joao Lucas Reis Da Silva', 'Olle Wallin', 'Gana Olle Wallin', 3.5),
    @('33', '14339505', '2025-08-04', 'Joel Schwaerzler', 'Christoph Negritu', 'Gana Christoph Negritu', 1.83),
    @('34', '14339509', '2025-08-04', 'Vilius Gaubas', 'Diego Dedura-Palomero', 'Gana Diego Dedura-Palomero', 2.63)
)

$firstRow = [int]$newRows[0][0]
$lastRow = [int]$newRows[$newRows.Count - 1][0]

# event_id (A) and fecha (B) are synced as plain text, matching the source feed
$ws.Range("A$firstRow" + ":B" + "$lastRow").NumberFormat = "@"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}
